# Viral Results terminology update (CVDLS-276)
#
# "Positive"      -> "Detected"
# "Negative"      -> "Not Detected"
# "Non-Negative"  -> "Inconclusive"
# "Recommended" remains unchanged (still a valid conclusion, just not new).
#
# These terms appear in the "Result" column of the Viral Results worksheet,
# so we do a whole-cell (not partial/substring) replace across the used
# cells of the sheet to make sure e.g. "Negative" doesn't also clobber
# "Non-Negative".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

$cells = $ws.UsedRange

# Order matters only in that each replace only matches whole-cell contents,
# so "Non-Negative" is never partially matched by the "Negative" replace.
$cells.Replace("Positive", "Detected", $whole) | Out-Null
$cells.Replace("Non-Negative", "Inconclusive", $whole) | Out-Null
$cells.Replace("Negative", "Not Detected", $whole) | Out-Null

# Move/record the active selection as it was left in the authored workbook.
$ws.Range("B14").Select()
